# Insert two new price records (rows 117-118) into the daily Cereza
# (cherry) price log, pushing all subsequent rows down by two.
#
# xlInsertShiftDown = -4121 (same semantics as Excel's Range.Insert with
# Shift:=xlShiftDown) — existing rows 117..215 move to 119..217 and the
# sheet's used range grows from A1:T215 to A1:T217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A117:T118").Insert(-4121)

# New row 117: Santina / Primera, Provincia de Curicó, 2023-12-20
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 45280
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100103
$ws.Range("H117").Value = "Frutos de hueso (carozo)"
$ws.Range("I117").Value = 100103001
$ws.Range("J117").Value = "Cereza"
$ws.Range("K117").Value = "Santina"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 100
$ws.Range("N117").Value = 9000
$ws.Range("O117").Value = 9000
$ws.Range("P117").Value = 9000
$ws.Range("Q117").Value = "$/bandeja 10 kilos"
$ws.Range("R117").Value = "Provincia de Curicó"
$ws.Range("S117").Value = 900
$ws.Range("T117").Value = 10

# New row 118: Santina / Segunda, Provincia de Curicó, 2023-12-20
$ws.Range("A118").Value = 7
$ws.Range("B118").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C118").Value = "Ñuble"
$ws.Range("D118").Value = 45280
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100103
$ws.Range("H118").Value = "Frutos de hueso (carozo)"
$ws.Range("I118").Value = 100103001
$ws.Range("J118").Value = "Cereza"
$ws.Range("K118").Value = "Santina"
$ws.Range("L118").Value = "Segunda"
$ws.Range("M118").Value = 100
$ws.Range("N118").Value = 7000
$ws.Range("O118").Value = 7000
$ws.Range("P118").Value = 7000
$ws.Range("Q118").Value = "$/bandeja 10 kilos"
$ws.Range("R118").Value = "Provincia de Curicó"
$ws.Range("S118").Value = 700
$ws.Range("T118").Value = 10
